$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vaccination")

# Remove the address-related columns (streetAndNr, zipCode, city, cantonCodeSender)
# from the first sheet, shifting the remaining columns left.
$ws.Range("E1:H1").EntireColumn.Delete()
